# Financials update: insert a new "most-recent period" column before column D
# on the PCB sheet, shifting the existing Period-Ending columns one slot to
# the right, then populate the new column with the latest period's figures.
# A few historical cells (rows 89, 94, 96) also get corrected figures for the
# periods that used to be in columns D/E (now E/F) as part of the same data
# refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at D; everything from D.. shifts right to E..
$ws.Columns("D").Insert()

# 2) The freshly inserted column D has no formatting of its own - copy the
#    number formats/styles from column E (the column that used to be D) so
#    the new column matches (date format in row 7/38/80, plain number style
#    elsewhere). Done per-block (skipping the fully blank separator rows 36
#    and 78) so no stray empty <row> entries get materialized there.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the new column D with the latest period's values.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 83700
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 19200
$ws.Range("D18").Value = 64500
$ws.Range("D20").Value = -29800
$ws.Range("D21").Value = 36000
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 34700
$ws.Range("D24").Value = 10400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 24300
$ws.Range("D27").Value = 24300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 29800
$ws.Range("D33").Value = 24300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 24300
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 24100
$ws.Range("D42").Value = 145600
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 4600
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 3400
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1697000
$ws.Range("D57").Value = 13000
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1486700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 37600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 210300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 24300
$ws.Range("D83").Value = 1300
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 35100
$ws.Range("E89").Value = 24700
$ws.Range("F89").Value = 40200
$ws.Range("D91").Value = -1100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -172500
$ws.Range("E94").Value = -219000
$ws.Range("F94").Value = -214000
$ws.Range("D96").Value = -1800
$ws.Range("E96").Value = -1600
$ws.Range("F96").Value = -1400
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 226100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 88600
